# Auto-generated edit script: update cryptos price/volume columns per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.333.96"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "'3.588.08"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'199.75"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").Value = "'592.05"
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E9").Value = "  +5.61%  "
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").Value = "'53.20"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").Value = "'0.0000300"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "'9.60"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "'696.08"
$ws.Range("E14").Value = "  +16.42%  "
$ws.Range("D15").Value = "'4.160.81"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "'70.396.08"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "'12.72"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "'19.02"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").Value = "'3.596.45"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "'0.991"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'18.26"
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("D23").Value = "'110.86"
$ws.Range("E23").Value = "  +7.88%  "
$ws.Range("D24").Value = "'5.30"
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("D25").Value = "'4.51"
$ws.Range("E25").Value = "  -2.93%  "
$ws.Range("D26").Value = "'2.99"
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("D27").Value = "'10.49"
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "'9.96"
$ws.Range("E29").Value = "  +3.89%  "
$ws.Range("D30").Value = "'34.56"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("D31").Value = "'4.40"
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("D33").Value = "'12.22"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").Value = "'63.53"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").Value = "'0.0₃0844"
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("D37").Value = "'3.800.22"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "'3.65"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").Value = "'510.12"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("D41").Value = "'2.99"
$ws.Range("E41").Value = "  -7.48%  "
$ws.Range("D42").Value = "'36.40"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("D45").Value = "'0.0469"
$ws.Range("E45").Value = "  +3.41%  "
$ws.Range("E46").Value = "  +6.92%  "
$ws.Range("E47").Value = "  +3.32%  "
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "'8.64"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "'1.81"
$ws.Range("E51").Value = "  +21.07%  "
